# repull data, push all data, mean calculation
# Update column F ("dSF") values for the rows that were recomputed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    6  = 0
    14 = -4
    19 = -11
    22 = -1
    28 = 4
    35 = -5
    39 = -6
    40 = -2
    44 = -2
    46 = 0
    50 = 0
    52 = -3
    54 = -2
    55 = 1
    59 = -2
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
